# feat: add 2022-Q1 data
#
# 1. The sheet that used to hold the "总计" (totals) table becomes the detail
#    sheet for the new "2022-Q1" quarter (it is renamed and repopulated with
#    the fund holding detail rows for 2022-Q1).
# 2. A brand-new "总计" sheet is appended at the end, containing the same
#    totals table as before plus a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: turn the old "总计" sheet into the "2022-Q1" detail sheet
# ---------------------------------------------------------------------
$detail = $wb.Worksheets.Item("总计")
$detail.Cells.Clear()
$detail.Name = "2022-Q1"

$detailHeaders = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 0; $col -lt $detailHeaders.Count; $col++) {
    $cell = $detail.Cells.Item(1, $col + 2)
    $cell.Value = $detailHeaders[$col]
    $cell.Style = "Bold"
}

$detailRows = @(
    @("009225", "天弘中证中美互联网指数（QDII）A", "1.84", "94.90", "10.20", "0.1877", 2),
    @("009226", "天弘中证中美互联网指数（QDII）C", "0.59", "94.90", "10.20", "0.0602", 2),
    @("009562", "工银瑞信中国机会全球配置股票(QDII)美元", "6.65", "92.85", "0.89", "0.0592", 10),
    @("486001", "工银瑞信中国机会全球配置股票(QDII)", "6.65", "92.85", "0.89", "0.0592", 10),
    @("009563", "工银瑞信中国机会全球配置股票(QDII)港币", "6.65", "92.85", "0.89", "0.0592", 10)
)

for ($i = 0; $i -lt $detailRows.Count; $i++) {
    $r = $i + 2
    $row = $detailRows[$i]

    $detail.Cells.Item($r, 1).Value = $i

    $detail.Cells.Item($r, 2).Value = "'" + $row[0]
    $detail.Cells.Item($r, 3).Value = $row[1]
    $detail.Cells.Item($r, 4).Value = "'" + $row[2]
    $detail.Cells.Item($r, 5).Value = "'" + $row[3]
    $detail.Cells.Item($r, 6).Value = "'" + $row[4]
    $detail.Cells.Item($r, 7).Value = "'" + $row[5]
    $detail.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# Step 2: add a fresh "总计" sheet after "2022-Q1" with the updated totals
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $detail)
$totals.Name = "总计"

$totalsHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($col = 0; $col -lt $totalsHeaders.Count; $col++) {
    $cell = $totals.Cells.Item(1, $col + 2)
    $cell.Value = $totalsHeaders[$col]
    $cell.Style = "Bold"
}

$totalsRows = @(
    @("2022-Q1", 5, 0.43),
    @("2021-Q4", 5, 0.46),
    @("2021-Q3", 5, 0.54),
    @("2021-Q2", 6, 1.25),
    @("2021-Q1", 5, 0.55),
    @("2020-Q4", 6, 0.6)
)

for ($i = 0; $i -lt $totalsRows.Count; $i++) {
    $r = $i + 2
    $row = $totalsRows[$i]

    $totals.Cells.Item($r, 1).Value = $i
    $totals.Cells.Item($r, 2).Value = $row[0]
    $totals.Cells.Item($r, 3).Value = $row[1]
    $totals.Cells.Item($r, 4).Value = $row[2]
}
